# DC-Colos.xlsx: insert a new "KCH" (Kuching, Malaysia) colo row just above
# the "IAD" (Ashburn) row, shifting every following row down by one.
# Net effect matches the target diff: row 275 becomes KCH/KUCHING, and all
# rows that used to be 275..335 become 276..336 (dimension grows to H336).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 275; this shifts rows 275-335 down to 276-336
# (values, formatting and styles all move together).
$ws.Rows.Item(275).Insert()

# The newly inserted row doesn't carry the bordered/bold "colo" style used by
# column A in every data row. Copy that formatting from the row directly
# below (the shifted-down "IAD" row) onto the new row before filling values.
$ws.Range("A276").Copy()
$ws.Range("A275").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new row's data.
$ws.Range("A275").Value2 = "KCH"
$ws.Range("B275").Value2 = "KUCHING, Malaysia"
$ws.Range("C275").Value2 = "Asia Pacific"
$ws.Range("D275").Value2 = "KUCHING"
$ws.Range("E275").Value2 = "Malaysia"
$ws.Range("F275").Value2 = "MY"
$ws.Range("G275").Value2 = 1.709727
$ws.Range("H275").Value2 = 110.353455
